# Chiffres COVID-19 Valais - data upload for the rows covering
# 2021-03-01 .. 2021-03-08 (rows 351 and 373-377).
#
# Column layout on this sheet:
#   A Date | B Cumul cas positifs (formula) | C Nb nouveaux cas positifs
#   D Nb nouvelles admissions a l'hopital | E Nb nouveaux deces a l'hopital
#   F Nb nouveaux deces extra-hospitaliers | G Patients COVID-19 aux SI total
#   H Total hospitalisations (formula) | I Nb de nouvelles sorties
#   J Cumul deces COVID-19 (formula) | K Nb nouveaux deces COVID-19 (formula)
#   L Nb nouveaux deces a l'hopital (raw) | M Nb nouveaux deces extra-hospitaliers (raw)
#
# B, H, J and K are pre-existing shared formulas that recompute automatically
# once the underlying raw cells below are filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correction to an already-reported day (2021-03-01): one more new case.
$ws.Range("C351").Value = 64

# Corrections to 2021-03-04 / 2021-03-05 new-case counts.
$ws.Range("C373").Value = 54
$ws.Range("C374").Value = 91

# New data for 2021-03-06 (row 375)
$ws.Range("C375").Value = 29
$ws.Range("E375").Value = 7
$ws.Range("F375").Value = 4
$ws.Range("G375").Value = 25
$ws.Range("L375").Value = 0
$ws.Range("M375").Value = 0

# New data for 2021-03-07 (row 376)
$ws.Range("C376").Value = 24
$ws.Range("E376").Value = 6
$ws.Range("F376").Value = 3
$ws.Range("G376").Value = 30
$ws.Range("L376").Value = 0
$ws.Range("M376").Value = 0

# New data for 2021-03-08 (row 377)
$ws.Range("C377").Value = 18
$ws.Range("E377").Value = 9
$ws.Range("F377").Value = 8
$ws.Range("G377").Value = 37
$ws.Range("L377").Value = 0
$ws.Range("M377").Value = 0
